$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the numeric estimate/SE/t.ratio columns (B, C, E) to 3 decimal places
# for data rows 2 through 28.
for ($r = 2; $r -le 28; $r++) {
    foreach ($col in @("B", "C", "E")) {
        $cell = $ws.Range("$col$r")
        $orig = [double]$cell.Value2
        $cell.Value = [Math]::Round($orig, 3)
    }
}
